# "Small fix on formatted agreement"
#
# 1) The {{tenant_addresss placeholder paragraph is missing its closing
#    "}}" -- every other {{...}} placeholder in the document is closed,
#    this one got dropped. Add a trailing run containing "}}" right
#    after the existing "tenant_addresss" run (after its spellEnd proof
#    mark), using the same Times New Roman / bCs run formatting already
#    used by the neighbouring placeholder runs in that paragraph.
#
# 2) Two built-in styles (Default Paragraph Font / Normal Table) pick up
#    an explicit "unhide when used" visibility flag.

$d = $word.ActiveDocument

# --- 1. Close the unterminated {{tenant_addresss placeholder -------------

$search = $d.Content
$found = $search.Find.Execute("tenant_addresss", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    $para = $search.Paragraphs(1).Range
    # Collapse to just before the paragraph mark (the last character of
    # a paragraph Range is always the pilcrow), then append the closing
    # braces there so they land in their own new run.
    $insertPoint = $d.Range($para.End - 1, $para.End - 1)
    $insertPoint.InsertAfter("}}")

    # Match the formatting used throughout the document for these
    # template-tag runs: Times New Roman in all four font slots, plus
    # "bold in complex scripts" (bCs) - no plain Bold.
    $insertPoint.Font.NameAscii = "Times New Roman"
    $insertPoint.Font.NameFarEast = "Times New Roman"
    $insertPoint.Font.NameOther = "Times New Roman"
    $insertPoint.Font.NameBi = "Times New Roman"
    $insertPoint.Font.BoldBi = $true
}

# --- 2. Style visibility flags -------------------------------------------

$tableNormal = $d.Styles("Normal Table")
$tableNormal.UnhideWhenUsed = $true

$defaultParaFont = $d.Styles("Default Paragraph Font")
$defaultParaFont.UnhideWhenUsed = $true
